$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ncam1"
$ws.Range("C2").Value = "Robo3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5869213333333333
$ws.Range("H2").Value = 1.760764
$ws.Range("I2").Value = 0.01834079054277507
$ws.Range("J2").Value = 0.02446997850397404
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.105595666666666
$ws.Range("N2").Value = 9.316787
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 1.822740349474222
$ws.Range("R2").Value = 16.404663145268
$ws.Range("S2").Value = 0.01834079054277507
$ws.Range("T2").Value = 0.02446997850397404

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Ncam1"
$ws.Range("C3").Value = "Robo3"
$ws.Range("D3").Value = "ECs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 7.268947999999999
$ws.Range("H3").Value = 21.806844
$ws.Range("I3").Value = 0.2271484186426865
$ws.Range("J3").Value = 0.3030576522007011
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.105595666666666
$ws.Range("N3").Value = 9.316787
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 22.57441341002533
$ws.Range("R3").Value = 203.169720690228
$ws.Range("S3").Value = 0.2271484186426865
$ws.Range("T3").Value = 0.3030576522007011

# Row 4
$ws.Range("A4").Value = "M1"
$ws.Range("B4").Value = "Ncam1"
$ws.Range("C4").Value = "Robo3"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.070339
$ws.Range("H4").Value = 0.211017
$ws.Range("I4").Value = 0.002198033693308568
$ws.Range("J4").Value = 0.002932580092490016
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.105595666666666
$ws.Range("N4").Value = 9.316787
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.2184444935976667
$ws.Range("R4").Value = 1.966000442379
$ws.Range("S4").Value = 0.002198033693308568
$ws.Range("T4").Value = 0.002932580092490016

# Row 5
$ws.Range("A5").Value = "M2"
$ws.Range("B5").Value = "Ncam1"
$ws.Range("C5").Value = "Robo3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.028134
$ws.Range("H5").Value = 0.084402
$ws.Range("I5").Value = 0.0008791634786895356
$ws.Range("J5").Value = 0.001172965329648049
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.105595666666666
$ws.Range("N5").Value = 9.316787
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.087372828486
$ws.Range("R5").Value = 0.7863554563740001
$ws.Range("S5").Value = 0.0008791634786895356
$ws.Range("T5").Value = 0.001172965329648049

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Ncam1"
$ws.Range("C6").Value = "Robo3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 24.046532
$ws.Range("H6").Value = 48.093064
$ws.Range("I6").Value = 0.7514335936425405
$ws.Range("J6").Value = 0.6683668238731867
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.105595666666666
$ws.Range("N6").Value = 9.316787
$ws.Range("O6").Value = 1
$ws.Range("P6").Value = 1
$ws.Range("Q6").Value = 74.67880557756132
$ws.Range("R6").Value = 448.072833465368
$ws.Range("S6").Value = 0.7514335936425405
$ws.Range("T6").Value = 0.6683668238731867

